$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F header: "Sunny/Rainy Test case", styled like the other
#     header cells (copy E1's format onto F1 so it reuses the same
#     "Accent6" header style rather than generating a brand-new one). ---
$ws.Range("F1").Value = "Sunny/Rainy Test case"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Shift the "Req"/"Design Mapping" (columns A/B) values of rows 6-18
#     up from the row below (the old "LTT_05" row for the Edit_data test
#     was effectively merged away), so every row from 6 down now carries
#     what used to be one row further down. ---
$ws.Range("A6").Value = "LTT_04"
$ws.Range("B6").Value = "7.1.4"

$ws.Range("A7").Value = "LTT_05"
$ws.Range("B7").Value = "7.1.5"

$ws.Range("A8").Value = "LTT_06"
$ws.Range("B8").Value = "7.1.6"

$ws.Range("A9").Value = "LTT_07"
$ws.Range("B9").Value = "7.1.7"

$ws.Range("A10").Value = "LTT_08"
$ws.Range("B10").Value = "7.1.8"

$ws.Range("A11").Value = "LTT_09"
$ws.Range("B11").Value = "7.1.9"

# From here down the "7.1.x" strings look like dates to Excel's
# smart-entry parser (e.g. "7.1.10" -> July 1 2010), so they're entered as
# a formula producing the literal text and then converted in place to a
# plain value (Copy + PasteSpecial values) - this avoids the auto
# date-conversion without leaving any extra number-format/style behind.
$ws.Range("A12").Value = "LTT_10"
$ws.Range("B12").Formula = '="7.1.10"'
$ws.Range("B12").Copy()
$ws.Range("B12").PasteSpecial(-4163)

$ws.Range("A13").Value = "LTT_11"
$ws.Range("B13").Formula = '="7.1.11"'
$ws.Range("B13").Copy()
$ws.Range("B13").PasteSpecial(-4163)

$ws.Range("A14").Value = "LTT_12"
$ws.Range("B14").Formula = '="7.1.12"'
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)

$ws.Range("A15").Value = "LTT_13"
$ws.Range("B15").Formula = '="7.1.13"'
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)

$ws.Range("A16").Value = "LTT_14"
$ws.Range("B16").Formula = '="7.1.14"'
$ws.Range("B16").Copy()
$ws.Range("B16").PasteSpecial(-4163)

$ws.Range("A17").Value = "LTT_15"
$ws.Range("B17").Formula = '="7.1.15"'
$ws.Range("B17").Copy()
$ws.Range("B17").PasteSpecial(-4163)

$ws.Range("A18").Value = "LTT_16"
$ws.Range("B18").Formula = '="7.1.16"'
$ws.Range("B18").Copy()
$ws.Range("B18").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- New column F ("Sunny"/"Rainy") values for the rows that carry a
#     UT Mapping (column D) test case but no IT Mapping (column E) entry.
#     The two report test-case rows (max_book_report/min_bal_report) are
#     "Rainy"; the rest are "Sunny". ---
$ws.Range("F4").Value = "Sunny"
$ws.Range("F6").Value = "Sunny"
$ws.Range("F7").Value = "Sunny"
$ws.Range("F8").Value = "Sunny"
$ws.Range("F10").Value = "Sunny"
$ws.Range("F11").Value = "Sunny"
$ws.Range("F12").Value = "Sunny"
$ws.Range("F13").Value = "Sunny"
$ws.Range("F15").Value = "Rainy"
$ws.Range("F16").Value = "Rainy"
$ws.Range("F18").Value = "Sunny"

# --- Column F width, matching the other data columns' widths (closest
#     achievable value given this engine's column-width quantisation) ---
$ws.Range("F1").ColumnWidth = 22.7

# --- Update the selection to mirror the saved workbook's last-active cell ---
$ws.Range("C22").Select()
